$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 415, shifting existing rows 415-527 down to 416-528
$ws.Rows(415).Insert()

# Populate the newly inserted row 415 with the new weekly data entry
$ws.Range("A415").Value = 3
$ws.Range("B415").Value = "Femacal de La Calera"
$ws.Range("C415").Value = "Coquimbo"
$ws.Range("D415").Value = 45204
$ws.Range("E415").Value = 5
$ws.Range("F415").Value = 100112001
$ws.Range("G415").Value = "Berenjena"
$ws.Range("H415").Value = "Sin especificar"
$ws.Range("I415").Value = "Primera"
$ws.Range("J415").Value = 50
$ws.Range("K415").Value = 9000
$ws.Range("L415").Value = 9000
$ws.Range("M415").Value = 9000
$ws.Range("N415").Value = "$/caja 60 unidades"
$ws.Range("O415").Value = "Región de Arica y Parinacota"
$ws.Range("P415").Value = 150
$ws.Range("Q415").Value = 60
$ws.Range("R415").Value = "Hortaliza"
